$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Numeric-looking price strings in column D must keep their exact original
# text formatting (e.g. "1.00", "165.40", "42.381.01"), so force the cell's
# number format to Text before assigning the value; this prevents Excel's
# COM layer from auto-converting the string into a numeric value and
# stripping formatting / dots.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.381.01"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.54"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.79"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.34"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.13"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.38"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.35"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.650.99"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.293.09"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.450.92"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.37"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "276.22"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +18.78%  "
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.83"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.40"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.89"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("E35").Value = "  -10.89%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0371"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.60"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("E39").Value = "  +3.65%  "
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.85"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.93"
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "82.78"
$ws.Range("E46").Value = "  +10.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.07"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.06"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.600.76"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.15"
$ws.Range("E51").Value = "  -3.49%  "
